$wb = $excel.ActiveWorkbook

# --- "Template" sheet: extend the trailing constant-XML cell so the
#     generated MODS record declares its type of resource before the
#     closing tags. Everything else on the row is unchanged; the shared
#     string table re-indexes itself as a consequence of this edit.
$ws1 = $wb.Worksheets.Item("Template")
$ws1.Range("DB1").Value = "<mods:typeOfResource>sound recording-nonmusical</mods:typeOfResource></mods:mods></datastream></object>"

# --- "Guide" sheet: document the new field with a row describing it,
#     mirroring the Field Label / description layout already used above.
$ws2 = $wb.Worksheets.Item("Guide")
$ws2.Range("A24").Value = "Type of Resource"
$ws2.Range("B24").Value = "constant value embedded in template- ""sound recording-nonmusical"""

# Leave the selection pointed at the newly-touched cells, matching where
# an author's cursor would land after typing these values in, with the
# Guide sheet left as the active/visible tab.
$ws1.Range("DB1").Select()
$ws2.Range("B24").Select()
